# Apply the "saving properties" edit to the Configs sheet:
# - Insert a new header/properties row above the existing "Year/Fica Type/Rate %/Rate" row
# - Add a new trailing blank row at the bottom of the table
# - Extend the box borders around the merged Year cells to cover the new bottom row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configs")

# 1. Insert a new row above row 3 (pushes the existing table down by one row)
$ws.Rows("3:3").Insert()

# 2. Populate the new row 3 ("properties" / "origin" / "Deviation" legend row)
$ws.Range("B3").Value = "properties"
$ws.Range("C3").Value = "origin"
$ws.Range("D3").Value = "Deviation"
$ws.Range("E3").Value = $null

# Font for the new legend row cells (matches the rest of the "Normal 2" styled cells)
$ws.Range("B3:E3").Font.Name = "Calibri"
$ws.Range("B3:E3").Font.Size = 11
$ws.Range("B3:E3").Font.Color = -16777216

# Thin border framing just the left edge of B3 and the right edge of E3
$ws.Range("B3").Borders.Item(7).LineStyle = 1
$ws.Range("B3").Borders.Item(7).Weight = 2
$ws.Range("E3").Borders.Item(10).LineStyle = 1
$ws.Range("E3").Borders.Item(10).Weight = 2

# 3. Add a new blank row 20 at the end of the table
$ws.Range("B20:E20").Value = $null
$ws.Range("B20:E20").Font.Name = "Calibri"
$ws.Range("B20:E20").Font.Size = 11
$ws.Range("B20:E20").Font.Color = -16777216

# 4. Re-apply the merged-cell box styling (font/fill/border) to the now 3-row-tall
#    Year groups so the newly shifted rows keep the same look as the visible cell.
$groups = @(5, 8, 11, 14, 17)
foreach ($top in $groups) {
  $r2 = $top + 1
  $r3 = $top + 2
  $rng = $ws.Range("B$r2" + ":B$r3")
  $rng.Font.Name = "Franklin Gothic Book"
  $rng.Font.Size = 10
  $rng.Font.Color = -16777215
  $rng.Interior.Color = 16448249
  $rng.HorizontalAlignment = -4108
  $rng.VerticalAlignment = -4108
  $rng.WrapText = $true
  $rng.Borders.Item(7).LineStyle = 1
  $rng.Borders.Item(7).Weight = 2
  $rng.Borders.Item(10).LineStyle = 1
  $rng.Borders.Item(10).Weight = 2
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(8).Weight = 2
  $rng.Borders.Item(9).LineStyle = 1
  $rng.Borders.Item(9).Weight = 2
}

# 5. Close the bottom of the table: the very last merged-group cell (B19) gets a
#    full thin box instead (closing border of the whole properties table).
$closeCell = $ws.Range("B19")
$closeCell.Borders.Item(7).LineStyle = 1
$closeCell.Borders.Item(7).Weight = 2
$closeCell.Borders.Item(10).LineStyle = 1
$closeCell.Borders.Item(10).Weight = 2
$closeCell.Borders.Item(8).LineStyle = 1
$closeCell.Borders.Item(8).Weight = 2
$closeCell.Borders.Item(9).LineStyle = 1
$closeCell.Borders.Item(9).Weight = 2
$closeCell.Font.Name = "Calibri"
$closeCell.Font.Size = 11
$closeCell.Font.Color = -16777216

Write-Output "done"
